$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.313941333333333
$ws.Range("H2").Value = 15.941824
$ws.Range("I2").Value = 0.176869630377001
$ws.Range("J2").Value = 0.176869630377001
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03057133333333333
$ws.Range("N2").Value = 0.09171399999999999
$ws.Range("O2").Value = 0.1084248755136686
$ws.Range("P2").Value = 0.1084248755136687
$ws.Range("Q2").Value = 0.1624542718151111
$ws.Range("R2").Value = 1.462088446336
$ws.Range("S2").Value = 0.01917706765577491
$ws.Range("T2").Value = 0.01917706765577492
$ws.Range("G3").Value = 5.313941333333333
$ws.Range("H3").Value = 15.941824
$ws.Range("I3").Value = 0.176869630377001
$ws.Range("J3").Value = 0.176869630377001
$ws.Range("O3").Value = 0.8915751244863314
$ws.Range("P3").Value = 0.8915751244863314
$ws.Range("Q3").Value = 1.335857541276444
$ws.Range("R3").Value = 12.022717871488
$ws.Range("S3").Value = 0.157692562721226
$ws.Range("T3").Value = 0.157692562721226
$ws.Range("I4").Value = 0.5461014638447835
$ws.Range("J4").Value = 0.5461014638447835
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03057133333333333
$ws.Range("N4").Value = 0.09171399999999999
$ws.Range("O4").Value = 0.1084248755136686
$ws.Range("P4").Value = 0.1084248755136687
$ws.Range("Q4").Value = 0.5015927011153333
$ws.Range("R4").Value = 4.514334310038
$ws.Range("S4").Value = 0.05921098323520287
$ws.Range("T4").Value = 0.05921098323520289
$ws.Range("I5").Value = 0.5461014638447835
$ws.Range("J5").Value = 0.5461014638447835
$ws.Range("O5").Value = 0.8915751244863314
$ws.Range("P5").Value = 0.8915751244863314
$ws.Range("S5").Value = 0.4868904806095807
$ws.Range("T5").Value = 0.4868904806095807
$ws.Range("I6").Value = 0.2770289057782155
$ws.Range("J6").Value = 0.2770289057782155
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03057133333333333
$ws.Range("N6").Value = 0.09171399999999999
$ws.Range("O6").Value = 0.1084248755136686
$ws.Range("P6").Value = 0.1084248755136687
$ws.Range("Q6").Value = 0.2544502923651111
$ws.Range("R6").Value = 2.290052631286
$ws.Range("S6").Value = 0.03003682462269086
$ws.Range("T6").Value = 0.03003682462269087
$ws.Range("I7").Value = 0.2770289057782155
$ws.Range("J7").Value = 0.2770289057782155
$ws.Range("O7").Value = 0.8915751244863314
$ws.Range("P7").Value = 0.8915751244863314
$ws.Range("S7").Value = 0.2469920811555247
$ws.Range("T7").Value = 0.2469920811555247
